$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for Case_0_3 (380 kV case): ext_grid vm_pu set-point changed from 1.05 to 1.02 pu,
# bus voltage magnitude results recomputed accordingly for rows 2-25, columns B-F and I-N.
$data = @{
    2 = @(1.02, 1.057425199387504, 1.063537018171536, 1.062994363749441, 1.072390465722183, 1.028591377447806, 1.062421689620653, 1.06625539010739, 1.065714204136541, 1.075085106733552, 1.024306450101033)
    3 = @(1.02, 1.059139243157481, 1.065086739482378, 1.064509015793468, 1.07396563762599, 1.028594865478385, 1.06378452191743, 1.067618298971018, 1.067042023848593, 1.076475145401226, 1.024781196286917)
    4 = @(1.02, 1.060246551893071, 1.066087892475023, 1.065487489949915, 1.074982864391474, 1.028594567344161, 1.064664250907067, 1.068498075854485, 1.067899101861595, 1.077372084601603, 1.025087135501995)
    5 = @(1.02, 1.060711646808538, 1.066508398562654, 1.065898464742028, 1.075410033333853, 1.028593829963519, 1.065033591438032, 1.068867435932898, 1.068258920549246, 1.077748565600155, 1.025215455204701)
    6 = @(1.02, 1.060789714127703, 1.066578981450272, 1.065967447443734, 1.075481729395417, 1.028593670253109, 1.065095576414467, 1.068929424151428, 1.06831930678998, 1.07781174392875, 1.025236983306637)
    7 = @(1.02, 1.060252768139541, 1.066093512772802, 1.065492982879125, 1.074988574090334, 1.028594559896501, 1.06466918799619, 1.068503013207361, 1.067903911717525, 1.077377117474561, 1.025088851279647)
    8 = @(1.02, 1.058004846463322, 1.064061094847477, 1.063506585660896, 1.072923225147435, 1.028593085098482, 1.062882708288527, 1.066716434877347, 1.066163388931907, 1.075555401319762, 1.024467154058374)
    9 = @(1.02, 1.054029483172698, 1.060466887248644, 1.059993588063504, 1.069267973101042, 1.028570931177275, 1.059718114857391, 1.063551663592841, 1.063079832215952, 1.072325713558483, 1.023361914289677)
    10 = @(1.02, 1.051368949510993, 1.058061537033398, 1.057642474583281, 1.066819905805024, 1.0285430309114, 1.057596681608843, 1.061430120880547, 1.061012488409351, 1.070158885610881, 1.022618361336561)
    11 = @(1.02, 1.050214307407181, 1.057017677333741, 1.056622127153141, 1.065757071181097, 1.028527838144724, 1.056675177927869, 1.06050857455068, 1.060114427686257, 1.069217255880808, 1.022294757731598)
    12 = @(1.02, 1.049785015355303, 1.056629580302067, 1.056242768105572, 1.065361854276611, 1.028521727507451, 1.056332442636275, 1.060165824349746, 1.059780403970926, 1.068866974639479, 1.022174306880257)
    13 = @(1.02, 1.049877118580464, 1.056712844956258, 1.056324158223907, 1.065446649412625, 1.028523059407197, 1.056405980956925, 1.060239365822785, 1.059852073447548, 1.068942134764822, 1.022200155371402)
    14 = @(1.02, 1.050178830367037, 1.056985604544096, 1.05659077655448, 1.065724411323504, 1.028527342570952, 1.056646856519118, 1.06048025188923, 1.060086826298144, 1.069188312152337, 1.022284806349033)
    15 = @(1.02, 1.050364670664108, 1.057153612540393, 1.056755001324087, 1.065895491912521, 1.028529919642658, 1.056795208242712, 1.060628610214948, 1.06023140613893, 1.069339921208573, 1.022336929364948)
    16 = @(1.02, 1.051445523056773, 1.05813076455982, 1.057710142286123, 1.06689038256886, 1.028543973656391, 1.057657776568101, 1.061491218788992, 1.061072028042838, 1.070221306400555, 1.022639802974345)
    17 = @(1.02, 1.052122803720871, 1.05874307511953, 1.058308653483534, 1.06751369216146, 1.028551956364952, 1.058198056497956, 1.062031525391502, 1.06159854779248, 1.070773264110488, 1.022829345889021)
    18 = @(1.02, 1.052517598858631, 1.059100001748218, 1.058657533797221, 1.067876988199641, 1.028556312295539, 1.058512912686776, 1.062346397581523, 1.06190537997532, 1.071094886394332, 1.022939745049706)
    19 = @(1.02, 1.052652171606764, 1.059221666894825, 1.058776455791733, 1.068000817321112, 1.028557746623221, 1.058620223348649, 1.062453713767994, 1.062009954994304, 1.071204496415807, 1.022977361617693)
    20 = @(1.02, 1.052050164000306, 1.058677403218655, 1.058244461850063, 1.067446844940279, 1.028551130950942, 1.058140118609184, 1.061973584593557, 1.061542086022423, 1.070714077995352, 1.022809026098558)
    21 = @(1.02, 1.050089995101464, 1.056905293703235, 1.056512274034207, 1.065642629402842, 1.028526094186577, 1.056575937117797, 1.060409329367169, 1.060017709773102, 1.069115833431743, 1.022259885678345)
    22 = @(1.02, 1.048855199870325, 1.05578900338318, 1.055421110174953, 1.064505737725501, 1.028507648590913, 1.055589877935013, 1.059423229327683, 1.05905670008124, 1.068107951234081, 1.021913170680008)
    23 = @(1.02, 1.04951001649576, 1.056380972403566, 1.055999756586743, 1.065108667223954, 1.028517683215221, 1.05611285626607, 1.059946228722676, 1.059566396905651, 1.06864253674586, 1.022097109395787)
    24 = @(1.02, 1.052082987539527, 1.058707078215344, 1.058273467969542, 1.067477051156232, 1.028551504847877, 1.058166299112307, 1.061999766409856, 1.061567599524648, 1.070740822658447, 1.022818208225519)
    25 = @(1.02, 1.055058970748221, 1.061397656321672, 1.060903346754979, 1.070214880295123, 1.028578974076938, 1.060538260329944, 1.06437185431008, 1.063879018378626, 1.073163041229467, 1.023648817187134)
}

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($r in $data.Keys) {
    $rowValues = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowValues[$i]
    }
}
